# Update cryptocurrency price/volume data to the latest scrape.
# Rows 18/19 and 48/49 also swap rank order (coin name, link, price, volume
# all move together), matching the source feed's updated ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.590.90"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.851.98"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.51"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5264"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3239"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.96"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7833"
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07765"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "1.839.98"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.73"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.039"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.00"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007949"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "26.613.86"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.640"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.490"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.013"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.80"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.171"
$ws.Range("E25").Value = "  -5.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.679"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.01"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "111.72"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.186"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.110"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08722"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04872"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7238"
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.133"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.876"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.110"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.271"
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01790"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4875"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9009"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.13"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.967"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.691"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4198"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.990"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.13"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1236"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8908"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.04"
$ws.Range("E51").Value = "  +1.64%  "
